$d = $word.ActiveDocument

$replacements = @(
    @{old="2023-11-14 Tuesday"; new="2023-11-15 Wednesday"},
    @{old="75×71="; new="65×88="},
    @{old="17×92="; new="59×33="},
    @{old="22×46="; new="68×11="},
    @{old="64×76="; new="88×86="},
    @{old="11×38="; new="88×39="},
    @{old="49×59="; new="64×84="},
    @{old="23×72="; new="62×52="},
    @{old="48×23="; new="39×50="},
    @{old="98×20="; new="84×49="},
    @{old="80×27="; new="32×87="},
    @{old="21×20="; new="32×51="},
    @{old="15×77="; new="70×50="},
    @{old="74×94="; new="52×55="},
    @{old="24×69="; new="25×42="},
    @{old="66×47="; new="47×27="},
    @{old="60×78="; new="87×52="},
    @{old="64×45="; new="63×35="},
    @{old="76×35="; new="81×22="},
    @{old="20×28="; new="45×34="},
    @{old="41×99="; new="40×46="},
    @{old="93×48="; new="24×15="},
    @{old="97×40="; new="43×47="},
    @{old="35×31="; new="17×93="},
    @{old="99×64="; new="72×65="},
    @{old="19×11="; new="31×83="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $true, $false, $false, $false, $true, 1, $false, $r.new, 2)
}

$d.Save()
